# Update the dated header line and the worked-example table to the new
# "output generated at c8c62b6" values.
$d = $word.ActiveDocument

# Header paragraph: date label.
$d.Paragraphs.Item(1).Range.Text = "2025-07-11 Friday"

# The table has 20 rows (5 data rows with 3 blank spacer rows after each);
# data lives in table rows 1, 5, 9, 13, 17 (1-based), 5 columns each.
$table = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Text = "35÷3=11, 2" },
    @{ Row = 1;  Col = 2; Text = "56÷9=6, 2" },
    @{ Row = 1;  Col = 3; Text = "50÷5=10, 0" },
    @{ Row = 1;  Col = 4; Text = "20÷8=2, 4" },
    @{ Row = 1;  Col = 5; Text = "68÷8=8, 4" },

    @{ Row = 5;  Col = 1; Text = "75÷5=15, 0" },
    @{ Row = 5;  Col = 2; Text = "26÷6=4, 2" },
    @{ Row = 5;  Col = 3; Text = "99÷8=12, 3" },
    @{ Row = 5;  Col = 4; Text = "94÷5=18, 4" },
    @{ Row = 5;  Col = 5; Text = "95÷8=11, 7" },

    @{ Row = 9;  Col = 1; Text = "48÷6=8, 0" },
    @{ Row = 9;  Col = 2; Text = "93÷7=13, 2" },
    @{ Row = 9;  Col = 3; Text = "36÷2=18, 0" },
    @{ Row = 9;  Col = 4; Text = "99÷9=11, 0" },
    @{ Row = 9;  Col = 5; Text = "19÷4=4, 3" },

    @{ Row = 13; Col = 1; Text = "96÷8=12, 0" },
    @{ Row = 13; Col = 2; Text = "40÷8=5, 0" },
    @{ Row = 13; Col = 3; Text = "95÷2=47, 1" },
    @{ Row = 13; Col = 4; Text = "24÷4=6, 0" },
    @{ Row = 13; Col = 5; Text = "52÷9=5, 7" },

    @{ Row = 17; Col = 1; Text = "84÷3=28, 0" },
    @{ Row = 17; Col = 2; Text = "97÷9=10, 7" },
    @{ Row = 17; Col = 3; Text = "83÷7=11, 6" },
    @{ Row = 17; Col = 4; Text = "18÷5=3, 3" },
    @{ Row = 17; Col = 5; Text = "19÷9=2, 1" }
)

foreach ($u in $updates) {
    $table.Cell($u.Row, $u.Col).Range.Text = $u.Text
}
